# "final commit of upload excel file"
# Updates a few contact-detail values and bumps the header/data row
# heights slightly, matching the final upload edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FirstName (B2: tintu -> rohan, B3: Maya -> mini)
$ws.Range("B2").Value = "rohan"
$ws.Range("B3").Value = "mini"

# Street (H3: dfbdf -> abcd)
$ws.Range("H3").Value = "abcd"

# Hobbies: trailing comma removed
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("L3").Value = "Reading ,Writing"

# Rows 1-3 got a touch taller (18.75 -> 19.5)
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
